$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.619088000000001
$ws.Range("H2").Value = 13.857264
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 33.6320075
$ws.Range("N2").Value = 67.264015
$ws.Range("O2").Value = 0.3908110491225105
$ws.Range("P2").Value = 0.3281482467916435
$ws.Range("Q2").Value = 155.34920225916
$ws.Range("R2").Value = 932.0952135549601
$ws.Range("S2").Value = 0.3908110491225105
$ws.Range("T2").Value = 0.3281482467916435

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.619088000000001
$ws.Range("H3").Value = 13.857264
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.015213
$ws.Range("N3").Value = 36.045639
$ws.Range("O3").Value = 0.1396193194224409
$ws.Range("P3").Value = 0.1758490515669409
$ws.Range("Q3").Value = 55.49932618574401
$ws.Range("R3").Value = 499.493935671696
$ws.Range("S3").Value = 0.1396193194224409
$ws.Range("T3").Value = 0.1758490515669409

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.619088000000001
$ws.Range("H4").Value = 13.857264
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.122567
$ws.Range("N4").Value = 0.367701
$ws.Range("O4").Value = 0.001424254495001488
$ws.Range("P4").Value = 0.001793833426290924
$ws.Range("Q4").Value = 0.566147758896
$ws.Range("R4").Value = 5.095329830064
$ws.Range("S4").Value = 0.001424254495001488
$ws.Range("T4").Value = 0.001793833426290924

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.619088000000001
$ws.Range("H5").Value = 13.857264
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1323796666666667
$ws.Range("N5").Value = 0.397139
$ws.Range("O5").Value = 0.001538279759615546
$ws.Range("P5").Value = 0.001937447037358482
$ws.Range("Q5").Value = 0.6114733297440001
$ws.Range("R5").Value = 5.503259967696001
$ws.Range("S5").Value = 0.001538279759615546
$ws.Range("T5").Value = 0.001937447037358482

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.619088000000001
$ws.Range("H6").Value = 13.857264
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 20.59650833333334
$ws.Range("N6").Value = 61.789525
$ws.Range("O6").Value = 0.2393357883858265
$ws.Range("P6").Value = 0.3014408863169769
$ws.Range("Q6").Value = 95.13708448440002
$ws.Range("R6").Value = 856.2337603596001
$ws.Range("S6").Value = 0.2393357883858265
$ws.Range("T6").Value = 0.3014408863169769

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.619088000000001
$ws.Range("H7").Value = 13.857264
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.558276
$ws.Range("N7").Value = 39.116552
$ws.Range("O7").Value = 0.2272713088146052
$ws.Range("P7").Value = 0.1908305348607893
$ws.Range("Q7").Value = 90.34139797228801
$ws.Range("R7").Value = 542.048387833728
$ws.Range("S7").Value = 0.2272713088146052
$ws.Range("T7").Value = 0.1908305348607893
